$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 30, pushing the existing rows 30-48 down to 31-49.
$ws.Rows(30).Insert()

# Populate the new row's values (this also grows the shared-string table).
$ws.Range("A30").Value = "@@@"
$ws.Range("B30").Value = "Find a better way to do Python documentation"
$ws.Range("C30").Value = "11/2-11/16"

# Match the formatting used by equivalent existing rows:
#  - A30 should use the same "text, quote-prefixed" style as A6 (the other "@@@" cell)
#  - B30 should use the same "indented detail" style as the other sprint-item B cells (e.g. B31)
$ws.Range("A6").Copy()
$ws.Range("A30").PasteSpecial(-4122)

$ws.Range("B31").Copy()
$ws.Range("B30").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Match the saved selection state from the target workbook.
$null = $ws.Range("B30").Select()
